$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "API Controllers" intro paragraph: the run-splitting collapses
#    into a single run. The visible text is unchanged, so a
#    find/replace over the exact text normalises it into one run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Developers will create the following API Controllers. For each controller, methods have been listed and descriptions given",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Developers will create the following API Controllers. For each controller, methods have been listed and descriptions given",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) & 3) The document carries exactly one "_GoBack" bookmark. Re-adding
#    it at the new location (the trailing empty paragraph, after the
#    "Hosting Recommendations" section is removed) automatically moves
#    it away from its old spot on the "Blazor Components" heading -
#    mirroring Word's own "only one _GoBack" behaviour. Do this BEFORE
#    any structural edits below, using a range that spans a paragraph
#    boundary (collapsed / single-paragraph ranges resolve unreliably
#    in this host).
# ---------------------------------------------------------------------
$docEnd = $d.Content.End
$bookmarkRange = $d.Range($docEnd - 2, $docEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# ---------------------------------------------------------------------
# Remove the standalone empty paragraph and the whole "Hosting
# Recommendations" section (heading, intro paragraph, "Web Server" /
# "Database" sub-headings and their placeholder paragraphs), leaving
# only the final (now-empty, now un-styled) paragraph that holds the
# relocated "_GoBack" bookmark.
# ---------------------------------------------------------------------
$startRng = $d.Content
$startRng.Find.Execute("Hosting Recommendations", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sectionStart = $startRng.Start - 1   # include the lone empty paragraph mark just before the heading

$endRng = $d.Content
$endRng.Find.Execute("Insert database server recommendations here.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sectionEnd = $endRng.End + 1         # include that paragraph's own end-of-paragraph mark

$d.Range($sectionStart, $sectionEnd).Delete() | Out-Null

# The trailing paragraph was styled "Heading 2"; strip that so it is a
# plain paragraph containing only the bookmark.
$d.Paragraphs.Last.Style = "Normal"
